# Auto-generated Excel COM-interop script applying the scheduled market-data
# refresh described by the commit "chore: update Sheets via scheduled runner".
# For each affected Leve row (identified by sheet + row number), update the
# currentAveragePrice / NQ / HQ and Leve price / profit columns (H:N) to their
# newly-fetched values. Columns M (LeveProfitNQ) and N (LeveProfitHQ) are only
# present when the corresponding HQ/NQ price path is viable, so some rows gain
# or lose one of those two cells entirely rather than merely changing value.

$wb = $excel.ActiveWorkbook

# ALC!row13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 1500
$ws.Range("I13").Value = 1500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -1331
$ws.Range("N13").ClearContents()

# ALC!row46
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 849
$ws.Range("I46").Value = 849
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2547
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2428
$ws.Range("N46").ClearContents()

# ALC!row60
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H60").Value = 849
$ws.Range("I60").Value = 849
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 2547
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -2063
$ws.Range("N60").ClearContents()

# ALC!row76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 14077.538
$ws.Range("I76").Value = 17022
$ws.Range("J76").Value = 12237.25
$ws.Range("K76").Value = 17022
$ws.Range("L76").Value = 12237.25
$ws.Range("M76").Value = -16707
$ws.Range("N76").Value = -12867.25

# ALC!row79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 14077.538
$ws.Range("I79").Value = 17022
$ws.Range("J79").Value = 12237.25
$ws.Range("K79").Value = 17022
$ws.Range("L79").Value = 12237.25
$ws.Range("M79").Value = -15930
$ws.Range("N79").Value = -14421.25

# ALC!row101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 198
$ws.Range("I101").Value = 198
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 594
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 1028

# ALC!row106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 22227790
$ws.Range("I106").Value = 23814774
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 23814774
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = -23814143
$ws.Range("N106").Value = -11262

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 11502.381
$ws.Range("I132").Value = 5505.2354
$ws.Range("J132").Value = 36990.25
$ws.Range("K132").Value = 16515.7062
$ws.Range("L132").Value = 110970.75
$ws.Range("M132").Value = -13985.7062
$ws.Range("N132").Value = -116030.75

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10427.866
$ws.Range("I45").Value = 19746.715
$ws.Range("J45").Value = 2273.875
$ws.Range("K45").Value = 19746.715
$ws.Range("L45").Value = 2273.875
$ws.Range("M45").Value = -19369.715
$ws.Range("N45").Value = -3027.875

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 22238.445
$ws.Range("I61").Value = 22878.143
$ws.Range("J61").Value = 19999.5
$ws.Range("K61").Value = 22878.143
$ws.Range("L61").Value = 19999.5
$ws.Range("M61").Value = -22666.143
$ws.Range("N61").Value = -20423.5

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5377.4
$ws.Range("I132").Value = 5329.3335
$ws.Range("J132").Value = 5449.5
$ws.Range("K132").Value = 15988.0005
$ws.Range("L132").Value = 16348.5
$ws.Range("M132").Value = -13458.0005
$ws.Range("N132").Value = -21408.5

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 22238.445
$ws.Range("I136").Value = 22878.143
$ws.Range("J136").Value = 19999.5
$ws.Range("K136").Value = 68634.429
$ws.Range("L136").Value = 59998.5
$ws.Range("M136").Value = -66084.429
$ws.Range("N136").Value = -65098.5

# BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 63111.5
$ws.Range("I22").Value = 698.8570999999999
$ws.Range("J22").Value = 500000
$ws.Range("K22").Value = 698.8570999999999
$ws.Range("L22").Value = 500000
$ws.Range("M22").Value = -525.8570999999999
$ws.Range("N22").Value = -500346

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3969.9546
$ws.Range("I105").Value = 2307.4211
$ws.Range("J105").Value = 14499.333
$ws.Range("K105").Value = 2307.4211
$ws.Range("L105").Value = 14499.333
$ws.Range("M105").Value = -560.4211
$ws.Range("N105").Value = -17993.333

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2976.52
$ws.Range("I134").Value = 2976.52
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 8929.559999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6394.559999999999

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2142.6428
$ws.Range("I31").Value = 1313.9
$ws.Range("J31").Value = 4214.5
$ws.Range("K31").Value = 1313.9
$ws.Range("L31").Value = 4214.5
$ws.Range("M31").Value = -1018.9
$ws.Range("N31").Value = -4804.5

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2142.6428
$ws.Range("I34").Value = 1313.9
$ws.Range("J34").Value = 4214.5
$ws.Range("K34").Value = 1313.9
$ws.Range("L34").Value = 4214.5
$ws.Range("M34").Value = -1111.9
$ws.Range("N34").Value = -4618.5

# CRP!row47
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 16298
$ws.Range("I47").Value = 16298
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 16298
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -15732

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10640.556
$ws.Range("I58").Value = 8852.75
$ws.Range("J58").Value = 14216.167
$ws.Range("K58").Value = 8852.75
$ws.Range("L58").Value = 14216.167
$ws.Range("M58").Value = -8649.75
$ws.Range("N58").Value = -14622.167

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10744.875
$ws.Range("I99").Value = 10778.444
$ws.Range("J99").Value = 10731.739
$ws.Range("K99").Value = 10778.444
$ws.Range("L99").Value = 10731.739
$ws.Range("M99").Value = -9280.444
$ws.Range("N99").Value = -13727.739

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10744.875
$ws.Range("I126").Value = 10778.444
$ws.Range("J126").Value = 10731.739
$ws.Range("K126").Value = 32335.332
$ws.Range("L126").Value = 32195.217
$ws.Range("M126").Value = -29865.332
$ws.Range("N126").Value = -37135.217

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9925
$ws.Range("I132").Value = 9928.4
$ws.Range("J132").Value = 9908
$ws.Range("K132").Value = 29785.2
$ws.Range("L132").Value = 29724
$ws.Range("M132").Value = -27255.2
$ws.Range("N132").Value = -34784

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11786.571
$ws.Range("I134").Value = 13383.818
$ws.Range("J134").Value = 5930
$ws.Range("K134").Value = 40151.454
$ws.Range("L134").Value = 17790
$ws.Range("M134").Value = -37616.454
$ws.Range("N134").Value = -22860

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 10640.556
$ws.Range("I136").Value = 8852.75
$ws.Range("J136").Value = 14216.167
$ws.Range("K136").Value = 26558.25
$ws.Range("L136").Value = 42648.501
$ws.Range("M136").Value = -24008.25
$ws.Range("N136").Value = -47748.501

# CUL!row17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 353.15384
$ws.Range("I17").Value = 256.85715
$ws.Range("J17").Value = 465.5
$ws.Range("K17").Value = 770.5714499999999
$ws.Range("L17").Value = 1396.5
$ws.Range("M17").Value = -601.5714499999999
$ws.Range("N17").Value = -1734.5

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 41767244
$ws.Range("I80").Value = 71002190
$ws.Range("J80").Value = 3035.2856
$ws.Range("K80").Value = 71002190
$ws.Range("L80").Value = 3035.2856
$ws.Range("M80").Value = -71001192
$ws.Range("N80").Value = -5031.2856

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 41767244
$ws.Range("I83").Value = 71002190
$ws.Range("J83").Value = 3035.2856
$ws.Range("K83").Value = 355010950
$ws.Range("L83").Value = 15176.428
$ws.Range("M83").Value = -355005958
$ws.Range("N83").Value = -25160.428

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 8311.362999999999
$ws.Range("I102").Value = 8851.875
$ws.Range("J102").Value = 6870
$ws.Range("K102").Value = 8851.875
$ws.Range("L102").Value = 6870
$ws.Range("M102").Value = -7229.875
$ws.Range("N102").Value = -10114

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3740.0535
$ws.Range("I126").Value = 3858.318
$ws.Range("J126").Value = 3663.5293
$ws.Range("K126").Value = 11574.954
$ws.Range("L126").Value = 10990.5879
$ws.Range("M126").Value = -9104.954000000002
$ws.Range("N126").Value = -15930.5879

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2399
$ws.Range("I132").Value = 2799
$ws.Range("J132").Value = 1599
$ws.Range("K132").Value = 8397
$ws.Range("L132").Value = 4797
$ws.Range("M132").Value = -5867
$ws.Range("N132").Value = -9857

# LTW!row11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2670.7222
$ws.Range("I122").Value = 2397.25
$ws.Range("J122").Value = 2748.8572
$ws.Range("K122").Value = 7191.75
$ws.Range("L122").Value = 8246.571599999999
$ws.Range("M122").Value = -4741.75
$ws.Range("N122").Value = -13146.5716

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18772.076
$ws.Range("I132").Value = 22085.277
$ws.Range("J132").Value = 3862.6667
$ws.Range("K132").Value = 66255.83099999999
$ws.Range("L132").Value = 11588.0001
$ws.Range("M132").Value = -63725.83099999999
$ws.Range("N132").Value = -16648.0001

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7519773
$ws.Range("I136").Value = 8196798
$ws.Range("J136").Value = 72498.5
$ws.Range("K136").Value = 24590394
$ws.Range("L136").Value = 217495.5
$ws.Range("M136").Value = -24587844
$ws.Range("N136").Value = -222595.5

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6952604.5
$ws.Range("I126").Value = 11368762
$ws.Range("J126").Value = 12927.857
$ws.Range("K126").Value = 34106286
$ws.Range("L126").Value = 38783.571
$ws.Range("M126").Value = -34103816
$ws.Range("N126").Value = -43723.571

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2174.4119
$ws.Range("I136").Value = 2312
$ws.Range("J136").Value = 1532.3334
$ws.Range("K136").Value = 6936
$ws.Range("L136").Value = 4597.0002
$ws.Range("N136").Value = -9697.0002

